{"js": "// Replace the generated two-digit-by-two-digit multiplication problems\n// (and their answers) throughout the document's table cells with a new\n// set of freshly generated problems, keeping everything else unchanged.\nconst replacements = [\n  [\"73\u00d769=5037\", \"55\u00d764=3520\"],\n  [\"66\u00d727=1782\", \"51\u00d744=2244\"],\n  [\"13\u00d728=364\", \"62\u00d743=2666\"],\n  [\"12\u00d749=588\", \"61\u00d721=1281\"],\n  [\"45\u00d784=3780\", \"99\u00d799=9801\"],\n  [\"58\u00d727=1566\", \"33\u00d772=2376\"],\n  [\"20\u00d716=320\", \"32\u00d795=3040\"],\n  [\"29\u00d797=2813\", \"95\u00d733=3135\"],\n  [\"76\u00d720=1520\", \"92\u00d742=3864\"],\n  [\"34\u00d719=646\", \"21\u00d766=1386\"],\n  [\"19\u00d714=266\", \"63\u00d722=1386\"],\n  [\"68\u00d718=1224\", \"99\u00d713=1287\"],\n  [\"95\u00d776=7220\", \"46\u00d711=506\"],\n  [\"38\u00d730=1140\", \"80\u00d728=2240\"],\n  [\"25\u00d744=1100\", \"13\u00d744=572\"],\n  [\"14\u00d727=378\", \"20\u00d745=900\"],\n  [\"70\u00d730=2100\", \"16\u00d728=448\"],\n  [\"25\u00d767=1675\", \"24\u00d741=984\"],\n  [\"70\u00d731=2170\", \"75\u00d712=900\"],\n  [\"90\u00d757=5130\", \"46\u00d758=2668\"],\n  [\"20\u00d798=1960\", \"24\u00d725=600\"],\n  [\"30\u00d793=2790\", \"49\u00d745=2205\"],\n  [\"99\u00d789=8811\", \"62\u00d738=2356\"],\n  [\"38\u00d768=2584\", \"87\u00d786=7482\"],\n  [\"97\u00d714=1358\", \"43\u00d725=1075\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the generated two-digit-by-two-digit multiplication problems\n# (and their answers) throughout the document's table cells with a new\n# set of freshly generated problems, keeping everything else unchanged.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"73\u00d769=5037\", \"55\u00d764=3520\"),\n  @(\"66\u00d727=1782\", \"51\u00d744=2244\"),\n  @(\"13\u00d728=364\",  \"62\u00d743=2666\"),\n  @(\"12\u00d749=588\",  \"61\u00d721=1281\"),\n  @(\"45\u00d784=3780\", \"99\u00d799=9801\"),\n  @(\"58\u00d727=1566\", \"33\u00d772=2376\"),\n  @(\"20\u00d716=320\",  \"32\u00d795=3040\"),\n  @(\"29\u00d797=2813\", \"95\u00d733=3135\"),\n  @(\"76\u00d720=1520\", \"92\u00d742=3864\"),\n  @(\"34\u00d719=646\",  \"21\u00d766=1386\"),\n  @(\"19\u00d714=266\",  \"63\u00d722=1386\"),\n  @(\"68\u00d718=1224\", \"99\u00d713=1287\"),\n  @(\"95\u00d776=7220\", \"46\u00d711=506\"),\n  @(\"38\u00d730=1140\", \"80\u00d728=2240\"),\n  @(\"25\u00d744=1100\", \"13\u00d744=572\"),\n  @(\"14\u00d727=378\",  \"20\u00d745=900\"),\n  @(\"70\u00d730=2100\", \"16\u00d728=448\"),\n  @(\"25\u00d767=1675\", \"24\u00d741=984\"),\n  @(\"70\u00d731=2170\", \"75\u00d712=900\"),\n  @(\"90\u00d757=5130\", \"46\u00d758=2668\"),\n  @(\"20\u00d798=1960\", \"24\u00d725=600\"),\n  @(\"30\u00d793=2790\", \"49\u00d745=2205\"),\n  @(\"99\u00d789=8811\", \"62\u00d738=2356\"),\n  @(\"38\u00d768=2584\", \"87\u00d786=7482\"),\n  @(\"97\u00d714=1358\", \"43\u00d725=1075\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
